$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read current data (rows 2-23: Language, Value)
$data = @()
for ($r = 2; $r -le 23; $r++) {
    $lang = $ws.Cells.Item($r, 1).Value()
    $val = $ws.Cells.Item($r, 2).Value()
    $data += [PSCustomObject]@{ Lang = $lang; Val = $val }
}

# Remove Swedish and Uzbek entries (dropped from the dataset)
$data = $data | Where-Object { $_.Lang -ne "Swedish" -and $_.Lang -ne "Uzbek" }

# Sort the remaining languages descending by their 2021 value
$data = $data | Sort-Object -Property Val -Descending

# Write sorted data back into rows 2-21, preserving existing cell formatting
$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row.Lang
    $ws.Cells.Item($r, 2).Value = $row.Val
    $r++
}

# Remove the now-unused trailing rows (22 and 23) so the sheet shrinks to A1:B21
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(22).Delete()
